$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for e880ddd3-... (row 3) moves from "In Translation" to "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 10:12:23"
$wsOverview.Range("E:F").Columns.AutoFit() | Out-Null

# --- zh-cn sheet: row for e880ddd3-... (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-26 10:12:19"
$wsZhCn.Range("C:C").Columns.AutoFit() | Out-Null

# --- de-de sheet: row for e880ddd3-... (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-26 10:12:23"
$wsDeDe.Range("C:C").Columns.AutoFit() | Out-Null
